$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 107.35229
$ws.Range("H2").Value = 322.05687
$ws.Range("I2").Value = 0.6580077109993711
$ws.Range("J2").Value = 0.6580077109993709
$ws.Range("M2").Value = 2.641449666666666
$ws.Range("N2").Value = 7.924348999999999
$ws.Range("O2").Value = 0.8024679156084781
$ws.Range("P2").Value = 0.8024679156084782
$ws.Range("Q2").Value = 283.5656706364033
$ws.Range("R2").Value = 2552.09103572763
$ws.Range("S2").Value = 0.5280300762999711
$ws.Range("T2").Value = 0.5280300762999711

# Row 3
$ws.Range("G3").Value = 107.35229
$ws.Range("H3").Value = 322.05687
$ws.Range("I3").Value = 0.6580077109993711
$ws.Range("J3").Value = 0.6580077109993709
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.650208
$ws.Range("N3").Value = 1.950624
$ws.Range("O3").Value = 0.1975320843915219
$ws.Range("P3").Value = 0.1975320843915219
$ws.Range("Q3").Value = 69.80131777632
$ws.Range("R3").Value = 628.2118599868801
$ws.Range("S3").Value = 0.1299776346993999
$ws.Range("T3").Value = 0.1299776346993999

# Row 4
$ws.Range("I4").Value = 0.1097031531157002
$ws.Range("J4").Value = 0.1097031531157002
$ws.Range("M4").Value = 2.641449666666666
$ws.Range("N4").Value = 7.924348999999999
$ws.Range("O4").Value = 0.8024679156084781
$ws.Range("P4").Value = 0.8024679156084782
$ws.Range("Q4").Value = 47.27611495150289
$ws.Range("R4").Value = 425.485034563526
$ws.Range("S4").Value = 0.08803326061643364
$ws.Range("T4").Value = 0.08803326061643363

# Row 5
$ws.Range("I5").Value = 0.1097031531157002
$ws.Range("J5").Value = 0.1097031531157002
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.650208
$ws.Range("N5").Value = 1.950624
$ws.Range("O5").Value = 0.1975320843915219
$ws.Range("P5").Value = 0.1975320843915219
$ws.Range("Q5").Value = 11.637287107264
$ws.Range("R5").Value = 104.735583965376
$ws.Range("S5").Value = 0.02166989249926654
$ws.Range("T5").Value = 0.02166989249926653

# Row 6
$ws.Range("G6").Value = 37.89738366666666
$ws.Range("H6").Value = 113.692151
$ws.Range("I6").Value = 0.2322891358849288
$ws.Range("J6").Value = 0.2322891358849288
$ws.Range("M6").Value = 2.641449666666666
$ws.Range("N6").Value = 7.924348999999999
$ws.Range("O6").Value = 0.8024679156084781
$ws.Range("P6").Value = 0.8024679156084782
$ws.Range("Q6").Value = 100.1040314538554
$ws.Range("R6").Value = 900.9362830846989
$ws.Range("S6").Value = 0.1864045786920734
$ws.Range("T6").Value = 0.1864045786920734

# Row 7
$ws.Range("G7").Value = 37.89738366666666
$ws.Range("H7").Value = 113.692151
$ws.Range("I7").Value = 0.2322891358849288
$ws.Range("J7").Value = 0.2322891358849288
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.650208
$ws.Range("N7").Value = 1.950624
$ws.Range("O7").Value = 0.1975320843915219
$ws.Range("P7").Value = 0.1975320843915219
$ws.Range("Q7").Value = 24.641182039136
$ws.Range("R7").Value = 221.770638352224
$ws.Range("S7").Value = 0.04588455719285545
$ws.Range("T7").Value = 0.04588455719285545
